# Update the "NegativeTests" class removal: corresponding test cases in the
# DataSet sheet's Status column (R) now reflect the new (passing) results.
# A leading apostrophe is used so Excel treats the assignment as literal text
# (matching the existing quote-prefixed text style of these cells) instead of
# re-evaluating/reformatting the cell, which would otherwise strip the
# quote-prefix styling already applied to these cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")

$passedRows = @(4, 5, 6, 7, 8, 10, 11, 12, 13, 18, 21, 23)
foreach ($row in $passedRows) {
    $ws.Range("R$row").Value = "'Passed"
}

$ws.Range("R25").Value = "'Inconclusive"
